# The post "「グレンダイザーのいる国は死にません」" (row 627) was removed from
# the source data. Delete its entire worksheet row so every following row
# shifts up by one (628->627, 629->628, ..., 670->669) and the sheet's used
# range shrinks from A1:C670 to A1:C669, matching a normal Excel row delete.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(627).Delete()
